# Fix typo on the "CONCLUSION" slide (slide 15): the text box that
# summarises the comparison between the considered algorithms.
#
# Before:
#   "Comparing the values of different considered algorithms we can
#    notice that the KRLS has the best accuracy."
# After:
#   "By comparing the values of the different algorithms considered we
#    can see that the KRLS has the best accuracy."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

# Work from the end of the paragraph towards the start so that the
# character offsets used below (taken from the original text) stay
# valid for every subsequent edit.

# "...we can notice..."  ->  "...we can see..."
$tr.Characters(64, 6).Text = "see"

# "...considered algorithms..." -> "...algorithms considered..."
# (same length swap, so it does not disturb any other offsets)
$tr.Characters(46, 10).Text = "considered"
$tr.Characters(35, 10).Text = "algorithms"

# "...values of different..." -> "...values of the different..."
$tr.Characters(21, 4).Text = " of the "

# "Comparing the values..." -> "By comparing the values..."
# First grow the leading run, then split it back into the two runs
# "By " and "comparing" so formatting matches the rest of the sentence.
$tr.Characters(1, 9).Text = "By Comparing"
$tr.Characters(4, 9).Text = "comparing"
